$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 3 (the "Engine | Investigate camera issues | 3" task),
# which shifts everything below it up by one row.
$ws.Rows("3:3").Delete()
